$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Properties")
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PrintQuality = @(0, 0)
$v = $ws.PageSetup.PrintQuality
Write-Host ("readback: " + $v)
